# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages have come back "in sync" with en-US:
#   - Overview sheet: per-language status text changes from
#     "Ready for handoff" to "Handed back: in sync with en-US".
#   - zh-cn / de-de detail sheets: the "Latest Target File" (I) and
#     "Latest Handback File" (J) columns get filled in (and the target-file
#     cell is turned into a hyperlink matching column A's), and the
#     "Latest Handback DateTime" (K) column gets a real timestamp instead of
#     the 0001-01-01 placeholder.
#   - A handful of columns are widened so the longer strings fit.

$wb = $excel.ActiveWorkbook

# Column-width quirk: this engine's Range/Columns.ColumnWidth (character
# units) gets translated into the raw OOXML <col width=.../> by adding a
# fixed 5/6 padding. Compute the character-width input that reproduces a
# desired raw width.
function RawWidthToColumnWidth($raw) {
    return $raw - (5.0 / 6.0)
}

# ----------------------------------------------------------------------
# Overview sheet (sheet 1): update zh-cn / de-de status text
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns (E, F) to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = RawWidthToColumnWidth(29.9777047293527)
$overview.Columns.Item(6).ColumnWidth = RawWidthToColumnWidth(29.9777047293527)

# ----------------------------------------------------------------------
# zh-cn sheet (sheet 2)
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Widen Status (C) and Latest Target File / Latest Handback File (I, J).
$zhcn.Columns.Item(3).ColumnWidth = RawWidthToColumnWidth(29.9777047293527)
$zhcn.Columns.Item(9).ColumnWidth = RawWidthToColumnWidth(40)
$zhcn.Columns.Item(10).ColumnWidth = RawWidthToColumnWidth(40)

# Fill in target file name + handback file name for both rows.
$zhcn.Range("I2").Value = "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md"
$zhcn.Range("J2").Value = "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.1aa5d01e0a11191b67d4d7b8421081020ab8e1fc.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-30 04:29:12"

$zhcn.Range("I3").Value = "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md"
$zhcn.Range("J3").Value = "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.ae7e940526978915f0e58db627d11bac240c7f26.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-30 04:29:12"

# Rebuild the sheet's hyperlinks in order (A2, I2, A3, I3) so the target
# file cells (I2/I3) pick up the same link + display text / style as the
# source file cells (A2/A3), and relationship ids renumber accordingly.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md", "", "", "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md", "", "", "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md", "", "", "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md", "", "", "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md") | Out-Null

# ----------------------------------------------------------------------
# de-de sheet (sheet 3)
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Widen Status (C) and Latest Target File / Latest Handback File (I, J).
$dede.Columns.Item(3).ColumnWidth = RawWidthToColumnWidth(29.9777047293527)
$dede.Columns.Item(9).ColumnWidth = RawWidthToColumnWidth(40)
$dede.Columns.Item(10).ColumnWidth = RawWidthToColumnWidth(40)

# Fill in target file name + handback file name for both rows.
$dede.Range("I2").Value = "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md"
$dede.Range("J2").Value = "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.1aa5d01e0a11191b67d4d7b8421081020ab8e1fc.de-de.xlf"
$dede.Range("K2").Value = "2016-08-30 04:29:21"

$dede.Range("I3").Value = "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md"
$dede.Range("J3").Value = "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.ae7e940526978915f0e58db627d11bac240c7f26.de-de.xlf"
$dede.Range("K3").Value = "2016-08-30 04:29:21"

# Rebuild the sheet's hyperlinks in order (A2, I2, A3, I3).
$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md", "", "", "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md", "", "", "7f51938f-be5a-46d3-b478-ca4d6c80a6ce.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md", "", "", "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/744621b1cfe592fb9e6e2b2e666cd515618c48a4/e2e/d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md", "", "", "d60885ef-e04d-4d3f-a8e8-3237afdb2a0e.md") | Out-Null

Write-Host "Handback report generated."
